# Clarify the REST interface note (cell A2 on the "REST" sheet) and
# update the row height / selection to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the explanatory note describing how write requests are handled.
$ws.Range("A2").Value = "The URL request shall be of the form http{s}://<host>/<top level>/<second level>/<third level>/. GET requests shall return the value of the item requested as a JSON object. In order to write an item, the user shall make a POST request with a JSON object containing the data to be written. The JSON object must contain a key-value pair with the key value equal to the second level name. All other keys shall be ignored. If the expected key is not present, return the current value. The only exception is the various waypoint writing commands, which expect a complete waypoint object. "

# The longer, wrapped text needs a taller row to display fully.
$ws.Rows.Item(2).RowHeight = 72

# Move the active selection / view to D6 (matches the saved view state).
$ws.Range("D6").Select()
